$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Change the "Runmode" column values from N to Y for rows 2, 4, 5, 6, 7
$ws.Range("C2").Value = "Y"
$ws.Range("C4").Value = "Y"
$ws.Range("C5").Value = "Y"
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"

# Update the selection shown in the sheet view
$ws.Range("C2:C7").Select()
